$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing values in row 74
$ws.Range("B74").Value = 965
$ws.Range("C74").Value = 468
$ws.Range("F74").Value = 4069
$ws.Range("L74").Value = 2220
$ws.Range("M74").Value = 1134

# Add new row 75 with the next quarter's data.
# A75 must become a shared-string text cell (like A74), not an auto-converted
# date serial. Enter it as a formula returning text, then paste-special as
# values so it collapses to a plain shared-string cell with no style change.
$ws.Range("A75").Formula = '="01-04-2021"'
$ws.Range("A75").Copy() | Out-Null
$ws.Range("A75").PasteSpecial(-4163) | Out-Null  # xlPasteValues

$ws.Range("B75").Value = -2842
$ws.Range("C75").Value = 726
$ws.Range("D75").Value = 61
$ws.Range("E75").Value = -3628
$ws.Range("F75").Value = -527
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 1
$ws.Range("I75").Value = -130
$ws.Range("J75").Value = 112
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = -146
$ws.Range("M75").Value = -311
$ws.Range("N75").Value = -53
